# Update the build timestamp embedded in the "version" strings across the
# workbook, going from "January 30 2026 16.19.47 EST" to
# "February 02 2026 12.49.33 EST".

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# --- "About" sheet ---------------------------------------------------
$aboutSheet = $wb.Worksheets.Item("About")

$aboutSheet.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"
$aboutSheet.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Karvina Coal Mines, Czech Republic, M0449, version ''mines - January 30 (built on ' + $newStamp + ')''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

# --- "Boundaries and methane sources" sheet ---------------------------
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 9; $row++) {
    $cell = $dataSheet.Range("S" + $row)
    $cell.Value = "mines - January 30 (built on " + $newStamp + ")"
}
